$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A14: Id
$ws.Range("A14").Value = 68946490

# C14: Valideringsstatus
$ws.Range("C14").Value = "Godkänd. Foto (eller ljud) granskat av validerare"

# I14: Antal (keep as text)
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "3"

# J14: Enhet
$ws.Range("J14").Value = "plantor/tuvor"

# K14: Ålder-Stadium
$ws.Range("K14").Value = "i frukt"

# L14: Kön - cell removed entirely
$ws.Range("L14").ClearContents()

# N14: Metod - cell removed entirely
$ws.Range("N14").ClearContents()

# P14: Lokalnamn
$ws.Range("P14").Value = "Skäftekärr NV-ut, Öl"

# S14: Noggrannhet
$ws.Range("S14").Value = 5

# Y14: Startdatum (keep as text, not an Excel date)
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = "2017-11-11"

# Z14: Starttid
$ws.Range("Z14").Value = "12:25"

# AA14: Slutdatum (keep as text, not an Excel date)
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = "2017-11-11"

# AB14: Sluttid
$ws.Range("AB14").Value = "13:12"

# AC14: Publik kommentar - new cell
$ws.Range("AC14").Value = "Tre ex, upp till 4 meter höga. Rikligt med bär."

# AF14: Bestämningsmetod - cell removed entirely
$ws.Range("AF14").ClearContents()

# AI14: Biotop-beskrivning - new cell
$ws.Range("AI14").Value = "Blandskog"

# AW14: Rapportör
$ws.Range("AW14").Value = "Joakim Ekman"

# AX14: Observatörer
$ws.Range("AX14").Value = "Joakim Ekman, Gabriel Ekman, Björn Owe-Larsson"
